$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Esdras (Ezra) row 16: mark "hecho" (done) column E as completed (1)
$ws.Range("E16").Value = 1

# Move active selection as recorded after the edit
$ws.Range("E27").Select()
